# Apply the changes described in the diff:
#  - F10: 85 -> 92
#  - F33: 75 -> 76
#  - Row 36 (Il Corrirere Della Sera / Facebook / 1 / 0) is removed entirely,
#    shifting the subsequent rows (37-39) up by one (they keep their own values).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update individual values
$ws.Range("F10").Value = 92
$ws.Range("F33").Value = 76

# Delete row 36 entirely, shifting rows below it up
$ws.Rows.Item(36).Delete()
